$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": add column BY (29-aug) ---
$ws1 = $wb.Worksheets.Item("Prix Spot")

# Copy header formatting from BX1 into BY1, then set its text
$ws1.Range("BX1").Copy($ws1.Range("BY1"))
$ws1.Range("BY1").Value = "29-aug"

# Fill the daily values for rows 2-25
$ws1.Range("BY2").Value = 82.7
$ws1.Range("BY3").Value = 51.1
$ws1.Range("BY4").Value = 51.93
$ws1.Range("BY5").Value = 25.04
$ws1.Range("BY6").Value = 17.87
$ws1.Range("BY7").Value = 21.33
$ws1.Range("BY8").Value = 30.57
$ws1.Range("BY9").Value = 45.47
$ws1.Range("BY10").Value = 57.53
$ws1.Range("BY11").Value = 42.47
$ws1.Range("BY12").Value = 5
$ws1.Range("BY13").Value = 0.83
$ws1.Range("BY14").Value = 3.08
$ws1.Range("BY15").Value = 0.05
$ws1.Range("BY16").Value = 0
$ws1.Range("BY17").Value = 0.62
$ws1.Range("BY18").Value = 5.66
$ws1.Range("BY19").Value = 3.42
$ws1.Range("BY20").Value = 22.36
$ws1.Range("BY21").Value = 53.55
$ws1.Range("BY22").Value = 81.8
$ws1.Range("BY23").Value = 63.5
$ws1.Range("BY24").Value = 94.39
$ws1.Range("BY25").Value = 86

# --- Sheet "Gaz": add row 74 ---
$ws2 = $wb.Worksheets.Item("Gaz")
$a74 = $ws2.Range("A74")
$a74.NumberFormat = "@"
$a74.Value = "2025-08-27"
$a74.Style = "Normal"
$ws2.Range("B74").Value = 31.4

# --- Sheet "CO2": add row 74 ---
$ws3 = $wb.Worksheets.Item("CO2")
$a74b = $ws3.Range("A74")
$a74b.NumberFormat = "@"
$a74b.Value = "2025-08-27"
$a74b.Style = "Normal"
$ws3.Range("B74").Value = 72.36
